$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new paragraph right after the title (Heading 1) that holds the
#    "Meta description" label (bold) followed by the description text.
# ---------------------------------------------------------------------------
$headingPara = $d.Paragraphs(2)
$headingPara.Range.InsertParagraphBefore()
$metaPara = $d.Paragraphs(2)

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Want to play Clone Bonus slot for free? Read our comprehensive review, ratings, and pros and cons. Learn how to win big and make your free play today!</w:t></w:r></w:p>'
$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2. Near the end of the document: drop the duplicated bold title paragraph
#    and turn the remaining italic paragraph into the new image-prompt text.
# ---------------------------------------------------------------------------
$dupTitlePara = $d.Paragraphs(55)
$dupTitlePara.Range.Delete()

$imgPara = $d.Paragraphs(55)
$imgXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Create a feature image fitting the game "Clone Bonus" that features a happy Maya warrior with glasses. The image should be in a cartoon style. The warrior should be standing in front of the slot machine, with a big smile on their face as they celebrate a big win. The background should be bright and colorful, with fruit symbols floating around in the air. The image should convey excitement and the potential for big wins, while also incorporating the Maya warrior theme of the game.</w:t></w:r></w:p>'
$imgPara.Range.InsertXML($imgXml)
